# Add a new "AddBooks" test case as row 7 of the APITestData sheet,
# mirroring the existing GetBooks/Register-style rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10
$json = "{" + $nl + "  ""userId"": ""string""," + $nl + "  ""collectionOfIsbns"": [" + $nl + "    {" + $nl + "      ""isbn"": ""9781449325862""" + $nl + "    }" + $nl + "  ]" + $nl + "}"

$ws.Range("A7").Value = "AddBooks"
$ws.Range("B7").Value = "/BookStore/v1/Books"
$ws.Range("C7").Value = "Authorization"
$ws.Range("D7").Value = "Bearer"
$ws.Range("G7").Value = $json
$ws.Range("G7").WrapText = $true

# Match the wrapped-text row height Excel auto-fits to for the 8-line body.
$ws.Rows.Item(7).RowHeight = 174

$ws.Range("G3").Select() | Out-Null
